$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.7
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 2.05
$ws.Range("P2").Value = 1.4
$ws.Range("Q2").Value = 2.75
$ws.Range("U2").Value = 19

# Row 4
$ws.Range("K4").Value = 19

# Row 5
$ws.Range("J5").Value = 1.14
$ws.Range("K5").Value = 5.5

# Row 7
$ws.Range("J7").Value = 1.04
$ws.Range("K7").Value = 13
$ws.Range("L7").Value = 1.22
$ws.Range("M7").Value = 4

# Row 8
$ws.Range("AD8").Value = 201
$ws.Range("AE8").Value = 21
$ws.Range("AG8").Value = 21
$ws.Range("AI8").Value = 41
$ws.Range("G8").Value = 1.42
$ws.Range("H8").Value = 4.5
$ws.Range("I8").Value = 7
$ws.Range("K8").Value = 15
$ws.Range("L8").Value = 1.17
$ws.Range("M8").Value = 5
$ws.Range("N8").Value = 1.53
$ws.Range("O8").Value = 2.4
$ws.Range("T8").Value = 9
$ws.Range("Z8").Value = 15

# Row 10
$ws.Range("AJ10").Value = 26
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 2.88
$ws.Range("J10").Value = 1.04
$ws.Range("K10").Value = 13
$ws.Range("L10").Value = 1.2
$ws.Range("M10").Value = 4.33
$ws.Range("N10").Value = 1.67
$ws.Range("O10").Value = 2.15
$ws.Range("P10").Value = 1.3
$ws.Range("Q10").Value = 3.4

# Row 14
$ws.Range("AB14").Value = 11.5
$ws.Range("AD14").Value = 350
$ws.Range("AE14").Value = 5.8
$ws.Range("AF14").Value = 8.25
$ws.Range("AG14").Value = 7.5
$ws.Range("AH14").Value = 16.5
$ws.Range("AI14").Value = 15
$ws.Range("G14").Value = 3.3
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 2.12
$ws.Range("N14").Value = 2
$ws.Range("O14").Value = 1.65
$ws.Range("P14").Value = 1.39
$ws.Range("T14").Value = 8
$ws.Range("U14").Value = 14.5
$ws.Range("V14").Value = 9.5
$ws.Range("W14").Value = 35
$ws.Range("X14").Value = 23
$ws.Range("Y14").Value = 28

# Row 17
$ws.Range("J17").Value = 1.08
$ws.Range("K17").Value = 8
$ws.Range("N17").Value = 2.25
$ws.Range("O17").Value = 1.62

# Row 19
$ws.Range("N19").Value = 1.73
$ws.Range("O19").Value = 2.08

# Row 20
$ws.Range("AB20").Value = 17
$ws.Range("AE20").Value = 6.2
$ws.Range("U20").Value = 16

# Row 21
$ws.Range("AE21").Value = 7.8
$ws.Range("AI21").Value = 32
$ws.Range("T21").Value = 6.8
$ws.Range("U21").Value = 10.75
$ws.Range("Z21").Value = 7.3

# Row 23
$ws.Range("AG23").Value = 9
$ws.Range("J23").Value = 1.06
$ws.Range("K23").Value = 10
$ws.Range("R23").Value = 1.8
$ws.Range("S23").Value = 1.91
$ws.Range("Z23").Value = 10

# Row 27
$ws.Range("G27").Value = 9
$ws.Range("I27").Value = 1.2
$ws.Range("J27").Value = 1.01
$ws.Range("K27").Value = 15
$ws.Range("Y27").Value = 51
$ws.Range("Z27").Value = 26

# Row 28
$ws.Range("AB28").Value = 15
$ws.Range("G28").Value = 7.5
$ws.Range("H28").Value = 5.5
$ws.Range("N28").Value = 1.22
$ws.Range("O28").Value = 4.2
$ws.Range("T28").Value = 41
$ws.Range("V28").Value = 26
$ws.Range("W28").Value = 101
$ws.Range("X28").Value = 51
$ws.Range("Y28").Value = 41

# Row 33
$ws.Range("N33").Value = 1.95
$ws.Range("O33").Value = 1.9

# Row 36
$ws.Range("AD36").Value = 251
$ws.Range("AE36").Value = 9
$ws.Range("AF36").Value = 15
$ws.Range("AG36").Value = 12
$ws.Range("AH36").Value = 34
$ws.Range("AI36").Value = 26
$ws.Range("G36").Value = 2.4
$ws.Range("I36").Value = 3
$ws.Range("L36").Value = 1.33
$ws.Range("M36").Value = 3.25
$ws.Range("N36").Value = 2.1
$ws.Range("O36").Value = 1.7
$ws.Range("T36").Value = 7.5
$ws.Range("U36").Value = 11
$ws.Range("V36").Value = 10
$ws.Range("W36").Value = 23
$ws.Range("X36").Value = 21
